# Applies the commit: re-sort applicant rows by final_score (desc),
# keep the top row's reasoning text (rewritten) and clear the
# reasoning/detail columns for the remaining rows (a "reasoning count limit").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (becomes id=2 / score=89.42 / name=Suk, with new reasoning text) ---
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = 89.42
$ws.Cells.Item(2, 3).Value = 'The applicant has a good understanding of web development technologies such as ReactJS, NodeJS, ExpressJS, and WebRTC, which are essential for the mentioned job role. The applicant has also worked on relevant projects showcasing skills in Flutter, Dart, Firebase, and Docker, demonstrating an ability to work with modern web and mobile technologies. While the applicant may not have direct experience with MongoDB, the applicant''s strong technical skills and adaptability make them well-suited to quickly learn and apply this skill in the given role.'
$ws.Cells.Item(2, 4).Value = '応募者は、ReactJS、NodeJS、ExpressJS、およびWebRTCなどのWeb開発技術について良い理解を持っており、これらは述べられた職務にとって不可欠です。応募者はまた、Flutter、Dart、Firebase、およびDockerのスキルを示す関連するプロジェクトで働いており、最新のWebおよびモバイル技術を扱う能力を示しています。応募者は直接的なMongoDBの経験を持っていないかもしれませんが、強力な技術スキルと適応性があり、与えられた役割で迅速にこのスキルを学び適用する能力を持っています。'
$ws.Cells.Item(2, 5).Value = 4
$ws.Cells.Item(2, 6).Value = 'The applicant has demonstrated a strong willingness to adapt to new cultures, as well as a keen interest in Japan''s work culture and language. Additionally, the applicant has showcased good teamwork and communication skills, along with a problem-solving approach. The applicant''s future career plans align with the company''s focus on AI/ML and backend development, indicating a good fit for the company''s goals.'
$ws.Cells.Item(2, 7).Value = '応募者は新しい文化に適応する意欲が強く、日本の労働文化や言語に強い興味を示しています。さらに、応募者は良いチームワークとコミュニケーションスキル、問題解決のアプローチを披露しています。応募者の将来のキャリアプランは、AI / MLおよびバックエンド開発に焦点を当てている会社の目標と一致しており、会社の目標に適しています。'
$ws.Cells.Item(2, 8).Value = "Suk"

# --- Row 3 (becomes id=4 / score=82.91 / name=Se, reasoning cleared) ---
$ws.Cells.Item(3, 1).Value = 4
$ws.Cells.Item(3, 2).Value = 82.91
$ws.Range("C3:G3").ClearContents()
$ws.Cells.Item(3, 8).Value = "Se"

# --- Row 4 (stays id=3 / score=77.25 / name=Ya, reasoning cleared) ---
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 77.25
$ws.Range("C4:G4").ClearContents()
$ws.Cells.Item(4, 8).Value = "Ya"

# --- Row 5 (becomes id=1 / score=74.02 / name=Shubh, reasoning cleared) ---
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = 74.02
$ws.Range("C5:G5").ClearContents()
$ws.Cells.Item(5, 8).Value = "Shubh"

Write-Output "applied reasoning count limit edit"
